$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match (bold, border, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for I2:J5
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9
